$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 9, shifting existing rows 9-25 down to 10-26.
$ws.Rows.Item(9).Insert()

# Populate the newly inserted row 9 with the new data record.
$ws.Cells.Item(9, 1).Value = 7
$ws.Cells.Item(9, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(9, 3).Value = "Ñuble"
$ws.Cells.Item(9, 4).Value = 44614
$ws.Cells.Item(9, 5).Value = 16
$ws.Cells.Item(9, 6).Value = 100112040
$ws.Cells.Item(9, 7).Value = "Cilantro"
$ws.Cells.Item(9, 8).Value = "Sin especificar"
$ws.Cells.Item(9, 9).Value = "Primera"
$ws.Cells.Item(9, 10).Value = 200
$ws.Cells.Item(9, 11).Value = 550
$ws.Cells.Item(9, 12).Value = 600
$ws.Cells.Item(9, 13).Value = 575
$ws.Cells.Item(9, 14).Value = "`$/atado 0,5 a 1 kilo"
$ws.Cells.Item(9, 15).Value = "Provincia de Diguillín"
$ws.Cells.Item(9, 16).Value = 575
$ws.Cells.Item(9, 17).Value = 1
$ws.Cells.Item(9, 18).Value = "Hortaliza"
